$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-9) in this sheet get their values permuted/updated per the
# new weekly price report. Only columns D, H, J, K, L, M, N, O, P, Q change;
# columns A, B, C, E, F, G, I, R stay identical across rows and are untouched.

$rows = @{
    2 = @{ D = 44438; H = "Española"; J = 400; K = 11000; L = 12000; M = 11500; N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 383; Q = 30 }
    3 = @{ D = 44498; H = "Española"; J = 400; K = 8500;  L = 9000;  M = 8750;  N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 292; Q = 30 }
    4 = @{ D = 44426; H = "Española"; J = 600; K = 11500; L = 12000; M = 11750; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 392; Q = 30 }
    5 = @{ D = 44426; H = "Madrigal"; J = 500; K = 12500; L = 13000; M = 12750; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 319; Q = 40 }
    6 = @{ D = 44427; H = "Madrigal"; J = 400; K = 12000; L = 13000; M = 12500; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 312; Q = 40 }
    7 = @{ D = 44420; H = "Madrigal"; J = 800; K = 14000; L = 15000; M = 14500; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 362; Q = 40 }
    8 = @{ D = 44420; H = "Madrigal"; J = 700; K = 13000; L = 14000; M = 13500; N = "`$/caja 40 unidades"; O = "Provincia del Elquí"; P = 338; Q = 40 }
    9 = @{ D = 44484; H = "Española"; J = 300; K = 9000;  L = 10000; M = 9500;  N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 317; Q = 30 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
}
